# LOQ4271.docx edit: move the PT/EN course-summary paragraphs down one slot
# (they now sit after "Objetivos"/"Programa resumido"/"Programa" rather than
# before them), relocate "Docente(s)" and "Bibliografia" to the end of the
# document, and rotate the three answer-runs inside the "Avaliação" bullet so
# the bibliography list becomes its final value run.
$d = $word.ActiveDocument

# --- Whole-paragraph content swaps. Done with Find/Replace (scoped to each
#     paragraphs own Range) rather than a bare Range.Text assignment so we
#     only touch the <w:t> text nodes and keep each paragraphs pPr/rPr intact
#     (also avoids Range.Text always stamping xml:space="preserve"). ---
$d.Paragraphs.Item(6).Range.Find.Execute("Objetivo Geral$([char]11)Permitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.$([char]11)$([char]11)Objetivos Específicos$([char]11)Saber planejar e executar um experimento fatorial completo e fracionado$([char]11)Saber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudada$([char]11)Dominar, pelo menos, um software comercial sobre o assunto$([char]11)Saber modelar um processo, com base em dados empíricos", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução Experimentação convencional Experimentos Fatoriais completos Experimentos Fatoriais fracionados Análise de variância Metodologia de superfície de resposta Método de Taguchi", 2) | Out-Null
$d.Paragraphs.Item(7).Range.Find.Execute("General objective To allow students to understand the mechanisms of obtaining the influence of several factors (independent variables of a process) on the response variables (dependent), through the multivariate analysis. Specific objectives Know how to plan and execute a complete and fractional factorial experiment Knowing to analyze the results proposing the condition of better fit that optimizes the values of the response variable in the studied experimental region Manage at least one commercial software on the subject Know how to model a process, based on empirical data", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction Conventional Experimentation, Full Factorial Experiments, Fractional Factorial Experiments, Analysis of Variance, Response Surface Methodology, Taguchi’s Method", 2) | Out-Null
$d.Paragraphs.Item(9).Range.Find.Execute("5840535 - Messias Borges Silva", $true, $false, $false, $false, $false, $true, 1, $false, "Objetivo Geral$([char]11)Permitir aos estudantes que compreendam os mecanismos de obtenção da influencia de diversos fatores (variáveis independentes de um processo) sobre as variáveis resposta (dependentes), através da análise multivariada.$([char]11)$([char]11)Objetivos Específicos$([char]11)Saber planejar e executar um experimento fatorial completo e fracionado$([char]11)Saber analisar os resultados propondo a condição de melhor ajuste que otimiza os valores da variável resposta na região experimental estudada$([char]11)Dominar, pelo menos, um software comercial sobre o assunto$([char]11)Saber modelar um processo, com base em dados empíricos", 2) | Out-Null
$d.Paragraphs.Item(11).Range.Find.Execute("Introdução Experimentação convencional Experimentos Fatoriais completos Experimentos Fatoriais fracionados Análise de variância Metodologia de superfície de resposta Método de Taguchi", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução Experimentação convencional Experimentos Fatoriais completos 2k , Experimentos Fatoriais fracionados 2k-p, Método de Plackett Burman,  Análise de variância Metodologia de superfície de resposta, Método de Taguchi .", 2) | Out-Null
$d.Paragraphs.Item(12).Range.Find.Execute("Introduction Conventional Experimentation, Full Factorial Experiments, Fractional Factorial Experiments, Analysis of Variance, Response Surface Methodology, Taguchi’s Method", $true, $false, $false, $false, $false, $true, 1, $false, "General objective To allow students to understand the mechanisms of obtaining the influence of several factors (independent variables of a process) on the response variables (dependent), through the multivariate analysis. Specific objectives Know how to plan and execute a complete and fractional factorial experiment Knowing to analyze the results proposing the condition of better fit that optimizes the values of the response variable in the studied experimental region Manage at least one commercial software on the subject Know how to model a process, based on empirical data", 2) | Out-Null
$d.Paragraphs.Item(14).Range.Find.Execute("Introdução Experimentação convencional Experimentos Fatoriais completos 2k , Experimentos Fatoriais fracionados 2k-p, Método de Plackett Burman,  Análise de variância Metodologia de superfície de resposta, Método de Taguchi .", $true, $false, $false, $false, $false, $true, 1, $false, "Provas, relatórios e apresentação de seminários.", 2) | Out-Null
$d.Paragraphs.Item(19).Range.Find.Execute("1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 1991$([char]11)2. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. $([char]11)3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. $([char]11)4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. $([char]11)5. COX, D.R., Planning of Experiments, Wiley 1976. $([char]11)6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. $([char]11)7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013", $true, $false, $false, $false, $false, $true, 1, $false, "5840535 - Messias Borges Silva", 2) | Out-Null

# --- Paragraph 17 (Avaliação bullet list) run-level rotation. Re-fetch the
#     paragraph Range before every Find call: Find.Execute collapses its Range
#     to the (newly replaced) match, so reusing a stale Range object would only
#     let us search the just-written replacement instead of the rest of the
#     paragraph. Processed slot "Norma -> Bibliografia" first, then "Critério ->
#     Norma-value", then "Método -> Critério-value" so every Find target is still
#     its untouched original (unique) text when it runs. ---
$d.Paragraphs.Item(17).Range.Find.Execute("Uma provas escrita com conteúdo de todo o semestre. NF = (MF + PR)/2, onde PR é uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "1. MONTGOMERY, D.C., Design and Analysis of Experiments, Wiley, 1991$([char]11)2. BOX, G.E.; HUNTER, W.G.; HUNTER, J.S., Statistic for Experimenters, John Wiley & Sons, New York, 1978. $([char]11)3. TAGUCHI, G.; WU, YU-IN., Introduction to off-Line Quality Control. Central Japan Quality Control Association. Meieki Nakamura-Ku Magaya, Japan, 1979. $([char]11)4. BRUNS, R.E., Como Fazer Experimentos, Editora UNICAMP, 2010. $([char]11)5. COX, D.R., Planning of Experiments, Wiley 1976. $([char]11)6. COX, G.M.; COCHRAN, W.G., Experimental Desing. Wiley 1976. $([char]11)7. SILVA M.B. et al, Design of Experiments-Applications, Editora Intech, 2013", 2) | Out-Null
$d.Paragraphs.Item(17).Range.Find.Execute("Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. MF = (0,40*Prova + 0,60*TRAB), onde TRAB é a nota média de trabalhos, relatórios e seminários.", $true, $false, $false, $false, $false, $true, 1, $false, "Uma provas escrita com conteúdo de todo o semestre. NF = (MF + PR)/2, onde PR é uma prova de recuperação.", 2) | Out-Null
$d.Paragraphs.Item(17).Range.Find.Execute("Provas, relatórios e apresentação de seminários.", $true, $false, $false, $false, $false, $true, 1, $false, "Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. MF = (0,40*Prova + 0,60*TRAB), onde TRAB é a nota média de trabalhos, relatórios e seminários.", 2) | Out-Null
